$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.233.82"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "3.139.93"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "635.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("E7").Value = "  +5.57%  "
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "3.139.64"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.727"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.03%  "
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "91.016.32"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "3.721.03"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "3.151.86"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("E19").Value = "  -4.07%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.91%  "
$ws.Range("E31").Value = "  -4.16%  "
$ws.Range("E32").Value = "  +11.24%  "
$ws.Range("E33").Value = "  +29.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "515.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.419"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0856"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +49.21%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.700"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.80%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.94%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("E51").Value = "  +3.29%  "
